# Updating graphs to be same structure (color, size, format, etc)
#
# 1) "NP 1" sheet: append a new "pDNA" summary block (rows 18-19) mirroring
#    the existing Average/STD block structure, using the shared string
#    "pDNA" as the row label.
# 2) "Sheet1": insert a new row 2 ("pDNA", N/P ratio 0) built from the new
#    NP1 pDNA summary values (halved average/stdev), shifting the existing
#    N/P-ratio rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "NP 1" worksheet — add the pDNA rows (18 & 19)
# ---------------------------------------------------------------------
$npws = $wb.Worksheets.Item("NP 1")

$npws.Range("A18").Value = "pDNA"
$npws.Range("B18").Value = 26.35
$npws.Range("C18").Value = 45.48
$npws.Range("D18").Value = 36.09
$npws.Range("E18").Formula = "=AVERAGE(B18:D18)"
$npws.Range("F18").Formula = "=STDEV(B18:D18)"

$npws.Range("E19").Formula = "=E18/2"
$npws.Range("F19").Formula = "=F18/2"

# ---------------------------------------------------------------------
# 2. "Sheet1" worksheet — insert a new row above row 2 for pDNA, pushing
#    the existing N/P-ratio rows down
# ---------------------------------------------------------------------
$sws = $wb.Worksheets.Item("Sheet1")

$sws.Rows.Item(2).Insert()
# Insert() copies the row-above's formatting into the newly inserted row;
# the authored pDNA row is unstyled, so strip that back off.
$sws.Range("A2:D2").ClearFormats()

$sws.Range("A2").Value = "pDNA"
$sws.Range("B2").Value = 0
$sws.Range("C2").Value = $npws.Range("E19").Value2
$sws.Range("D2").Value = $npws.Range("F19").Value2

# Selection / view bookkeeping to mirror the authored workbook state
# (Sheet1 stays the active/visible tab, so select it last.)
$npws.Range("E19:F19").Select()
$sws.Range("G6").Select()
